# Insert a new weekly price-report row for "Femacal de La Calera / Albahaca"
# at spreadsheet row 39, pushing all the existing rows below it down by one
# (old row 39 -> new row 40, ..., old row 152 -> new row 153).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 39, shifting rows 39:152 down to 40:153.
$ws.Rows(39).Insert()

# Populate the newly inserted row 39 with the new record's data.
$ws.Cells.Item(39, 1).Value  = 3
$ws.Cells.Item(39, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(39, 3).Value  = "Coquimbo"
$ws.Cells.Item(39, 4).Value  = 44624
$ws.Cells.Item(39, 5).Value  = 5
$ws.Cells.Item(39, 6).Value  = 100112052
$ws.Cells.Item(39, 7).Value  = "Albahaca"
$ws.Cells.Item(39, 8).Value  = "Sin especificar"
$ws.Cells.Item(39, 9).Value  = "Primera"
$ws.Cells.Item(39, 10).Value = 110
$ws.Cells.Item(39, 11).Value = 4500
$ws.Cells.Item(39, 12).Value = 5000
$ws.Cells.Item(39, 13).Value = 4727
$ws.Cells.Item(39, 14).Value = "`$/docena de matas"
$ws.Cells.Item(39, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(39, 16).Value = 788
$ws.Cells.Item(39, 17).Value = 6
$ws.Cells.Item(39, 18).Value = "Hortaliza"
